# Auto-generated COM-interop script implementing the commit:
#   feat: add 2022-Q1 data
# - inserts a new '2022-Q1' worksheet (fund holdings detail) between
#   '2021-Q4' and the summary sheet
# - rebuilds the '总计' (summary) worksheet with a new leading
#   '2022-Q1' row

$wb = $excel.ActiveWorkbook

# ---- locate anchor sheets before we start mutating the tab order ----
$q4 = $wb.Worksheets.Item('2021-Q4')
$styleSrcHeader = $wb.Worksheets.Item('2020-Q4').Range('B1')
$styleSrcIndex = $wb.Worksheets.Item('2020-Q4').Range('A2')

# ---- drop the old summary sheet; we rebuild it from scratch at the end ----
$wb.Worksheets.Item('总计').Delete()

# ---- new '2022-Q1' detail sheet, inserted right after '2021-Q4' ----
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = '2022-Q1'

# header row
$q1.Range('B1').Value = '基金代码'
$q1.Range('C1').Value = '基金名称'
$q1.Range('D1').Value = '基金规模'
$q1.Range('E1').Value = '股票总仓位'
$q1.Range('F1').Value = '仓位占比'
$q1.Range('G1').Value = '持有市值(亿元)'
$q1.Range('H1').Value = '仓位排名'

# text-formatted columns (preserve fund codes / decimal text exactly,
# e.g. leading zeros in '001044' and trailing zeros in '0.00')
$q1.Range('B2:B31').NumberFormat = '@'
$q1.Range('D2:D31').NumberFormat = '@'
$q1.Range('E2:E31').NumberFormat = '@'
$q1.Range('F2:F31').NumberFormat = '@'
$q1.Range('G2:G31').NumberFormat = '@'

# data rows
$q1.Cells.Item(2, 1).Value = 0
$q1.Range('B2').Value = '001044'
$q1.Range('C2').Value = '嘉实新消费股票'
$q1.Range('D2').Value = '8.92'
$q1.Range('E2').Value = '80.25'
$q1.Range('F2').Value = '4.63'
$q1.Range('G2').Value = '0.4130'
$q1.Cells.Item(2, 8).Value = 10
$q1.Cells.Item(3, 1).Value = 1
$q1.Range('B3').Value = '003713'
$q1.Range('C3').Value = '英大睿盛灵活配置混合A'
$q1.Range('D3').Value = '5.99'
$q1.Range('E3').Value = '87.42'
$q1.Range('F3').Value = '6.02'
$q1.Range('G3').Value = '0.3606'
$q1.Cells.Item(3, 8).Value = 5
$q1.Cells.Item(4, 1).Value = 2
$q1.Range('B4').Value = '005358'
$q1.Range('C4').Value = '东方阿尔法精选灵活配置混合A'
$q1.Range('D4').Value = '4.31'
$q1.Range('E4').Value = '92.94'
$q1.Range('F4').Value = '7.79'
$q1.Range('G4').Value = '0.3357'
$q1.Cells.Item(4, 8).Value = 6
$q1.Cells.Item(5, 1).Value = 3
$q1.Range('B5').Value = '167508'
$q1.Range('C5').Value = '安信价值发现两年定期开放混合（LOF）'
$q1.Range('D5').Value = '5.43'
$q1.Range('E5').Value = '89.42'
$q1.Range('F5').Value = '5.11'
$q1.Range('G5').Value = '0.2775'
$q1.Cells.Item(5, 8).Value = 6
$q1.Cells.Item(6, 1).Value = 4
$q1.Range('B6').Value = '001305'
$q1.Range('C6').Value = '九泰天富改革新动力混合A'
$q1.Range('D6').Value = '3.74'
$q1.Range('E6').Value = '88.86'
$q1.Range('F6').Value = '6.05'
$q1.Range('G6').Value = '0.2263'
$q1.Cells.Item(6, 8).Value = 7
$q1.Cells.Item(7, 1).Value = 5
$q1.Range('B7').Value = '003714'
$q1.Range('C7').Value = '英大睿盛灵活配置混合C'
$q1.Range('D7').Value = '2.40'
$q1.Range('E7').Value = '87.42'
$q1.Range('F7').Value = '6.02'
$q1.Range('G7').Value = '0.1445'
$q1.Cells.Item(7, 8).Value = 5
$q1.Cells.Item(8, 1).Value = 6
$q1.Range('B8').Value = '001782'
$q1.Range('C8').Value = '九泰久益灵活配置混合A'
$q1.Range('D8').Value = '2.33'
$q1.Range('E8').Value = '94.33'
$q1.Range('F8').Value = '5.94'
$q1.Range('G8').Value = '0.1384'
$q1.Cells.Item(8, 8).Value = 7
$q1.Cells.Item(9, 1).Value = 7
$q1.Range('B9').Value = '000520'
$q1.Range('C9').Value = '上银新兴价值成长混合'
$q1.Range('D9').Value = '4.64'
$q1.Range('E9').Value = '77.73'
$q1.Range('F9').Value = '2.61'
$q1.Range('G9').Value = '0.1211'
$q1.Cells.Item(9, 8).Value = 8
$q1.Cells.Item(10, 1).Value = 8
$q1.Range('B10').Value = '001399'
$q1.Range('C10').Value = '安信鑫安得利灵活配置混合A'
$q1.Range('D10').Value = '7.07'
$q1.Range('E10').Value = '23.18'
$q1.Range('F10').Value = '1.57'
$q1.Range('G10').Value = '0.1110'
$q1.Cells.Item(10, 8).Value = 5
$q1.Cells.Item(11, 1).Value = 9
$q1.Range('B11').Value = '004138'
$q1.Range('C11').Value = '上银鑫达灵活配置混合'
$q1.Range('D11').Value = '2.41'
$q1.Range('E11').Value = '77.12'
$q1.Range('F11').Value = '3.89'
$q1.Range('G11').Value = '0.0937'
$q1.Cells.Item(11, 8).Value = 7
$q1.Cells.Item(12, 1).Value = 10
$q1.Range('B12').Value = '009766'
$q1.Range('C12').Value = '安信平稳双利3个月持有期混合A'
$q1.Range('D12').Value = '2.33'
$q1.Range('E12').Value = '39.45'
$q1.Range('F12').Value = '4.01'
$q1.Range('G12').Value = '0.0934'
$q1.Cells.Item(12, 8).Value = 4
$q1.Cells.Item(13, 1).Value = 11
$q1.Range('B13').Value = '003345'
$q1.Range('C13').Value = '安信新成长灵活配置混合A'
$q1.Range('D13').Value = '6.27'
$q1.Range('E13').Value = '30.76'
$q1.Range('F13').Value = '1.48'
$q1.Range('G13').Value = '0.0928'
$q1.Cells.Item(13, 8).Value = 4
$q1.Cells.Item(14, 1).Value = 12
$q1.Range('B14').Value = '001844'
$q1.Range('C14').Value = '九泰久益灵活配置混合C'
$q1.Range('D14').Value = '1.47'
$q1.Range('E14').Value = '94.33'
$q1.Range('F14').Value = '5.94'
$q1.Range('G14').Value = '0.0873'
$q1.Cells.Item(14, 8).Value = 7
$q1.Cells.Item(15, 1).Value = 13
$q1.Range('B15').Value = '080001'
$q1.Range('C15').Value = '长盛成长价值混合'
$q1.Range('D15').Value = '2.53'
$q1.Range('E15').Value = '60.74'
$q1.Range('F15').Value = '2.23'
$q1.Range('G15').Value = '0.0564'
$q1.Cells.Item(15, 8).Value = 7
$q1.Cells.Item(16, 1).Value = 14
$q1.Range('B16').Value = '005359'
$q1.Range('C16').Value = '东方阿尔法精选灵活配置混合C'
$q1.Range('D16').Value = '0.54'
$q1.Range('E16').Value = '92.94'
$q1.Range('F16').Value = '7.79'
$q1.Range('G16').Value = '0.0421'
$q1.Cells.Item(16, 8).Value = 6
$q1.Cells.Item(17, 1).Value = 15
$q1.Range('B17').Value = '003446'
$q1.Range('C17').Value = '英大睿鑫灵活配置混合A'
$q1.Range('D17').Value = '0.59'
$q1.Range('E17').Value = '89.46'
$q1.Range('F17').Value = '7.03'
$q1.Range('G17').Value = '0.0415'
$q1.Cells.Item(17, 8).Value = 2
$q1.Cells.Item(18, 1).Value = 16
$q1.Range('B18').Value = '003447'
$q1.Range('C18').Value = '英大睿鑫灵活配置混合C'
$q1.Range('D18').Value = '0.51'
$q1.Range('E18').Value = '89.46'
$q1.Range('F18').Value = '7.03'
$q1.Range('G18').Value = '0.0359'
$q1.Cells.Item(18, 8).Value = 2
$q1.Cells.Item(19, 1).Value = 17
$q1.Range('B19').Value = '009912'
$q1.Range('C19').Value = '九泰天富改革新动力混合C'
$q1.Range('D19').Value = '0.59'
$q1.Range('E19').Value = '88.86'
$q1.Range('F19').Value = '6.05'
$q1.Range('G19').Value = '0.0357'
$q1.Cells.Item(19, 8).Value = 7
$q1.Cells.Item(20, 1).Value = 18
$q1.Range('B20').Value = '001400'
$q1.Range('C20').Value = '安信鑫安得利灵活配置混合C'
$q1.Range('D20').Value = '2.21'
$q1.Range('E20').Value = '23.18'
$q1.Range('F20').Value = '1.57'
$q1.Range('G20').Value = '0.0347'
$q1.Cells.Item(20, 8).Value = 5
$q1.Cells.Item(21, 1).Value = 19
$q1.Range('B21').Value = '007393'
$q1.Range('C21').Value = '上银未来生活灵活配置混合'
$q1.Range('D21').Value = '1.05'
$q1.Range('E21').Value = '85.99'
$q1.Range('F21').Value = '3.03'
$q1.Range('G21').Value = '0.0318'
$q1.Cells.Item(21, 8).Value = 8
$q1.Cells.Item(22, 1).Value = 20
$q1.Range('B22').Value = '004249'
$q1.Range('C22').Value = '安信中国制造2025沪港深灵活配置混合'
$q1.Range('D22').Value = '0.58'
$q1.Range('E22').Value = '89.89'
$q1.Range('F22').Value = '5.34'
$q1.Range('G22').Value = '0.0310'
$q1.Cells.Item(22, 8).Value = 5
$q1.Cells.Item(23, 1).Value = 21
$q1.Range('B23').Value = '004393'
$q1.Range('C23').Value = '安信合作创新主题沪港深灵活配置混合'
$q1.Range('D23').Value = '0.49'
$q1.Range('E23').Value = '89.26'
$q1.Range('F23').Value = '5.82'
$q1.Range('G23').Value = '0.0285'
$q1.Cells.Item(23, 8).Value = 6
$q1.Cells.Item(24, 1).Value = 22
$q1.Range('B24').Value = '009899'
$q1.Range('C24').Value = '上银内需增长股票'
$q1.Range('D24').Value = '0.59'
$q1.Range('E24').Value = '89.26'
$q1.Range('F24').Value = '3.10'
$q1.Range('G24').Value = '0.0183'
$q1.Cells.Item(24, 8).Value = 9
$q1.Cells.Item(25, 1).Value = 23
$q1.Range('B25').Value = '003346'
$q1.Range('C25').Value = '安信新成长灵活配置混合C'
$q1.Range('D25').Value = '1.18'
$q1.Range('E25').Value = '30.76'
$q1.Range('F25').Value = '1.48'
$q1.Range('G25').Value = '0.0175'
$q1.Cells.Item(25, 8).Value = 4
$q1.Cells.Item(26, 1).Value = 24
$q1.Range('B26').Value = '009767'
$q1.Range('C26').Value = '安信平稳双利3个月持有期混合C'
$q1.Range('D26').Value = '0.26'
$q1.Range('E26').Value = '39.45'
$q1.Range('F26').Value = '4.01'
$q1.Range('G26').Value = '0.0104'
$q1.Cells.Item(26, 8).Value = 4
$q1.Cells.Item(27, 1).Value = 25
$q1.Range('B27').Value = '008443'
$q1.Range('C27').Value = '九泰动态策略灵活配置混合A'
$q1.Range('D27').Value = '0.24'
$q1.Range('E27').Value = '64.82'
$q1.Range('F27').Value = '2.78'
$q1.Range('G27').Value = '0.0067'
$q1.Cells.Item(27, 8).Value = 9
$q1.Cells.Item(28, 1).Value = 26
$q1.Range('B28').Value = '750005'
$q1.Range('C28').Value = '安信平稳增长混合A'
$q1.Range('D28').Value = '0.08'
$q1.Range('E28').Value = '65.16'
$q1.Range('F28').Value = '5.39'
$q1.Range('G28').Value = '0.0043'
$q1.Cells.Item(28, 8).Value = 4
$q1.Cells.Item(29, 1).Value = 27
$q1.Range('B29').Value = '008444'
$q1.Range('C29').Value = '九泰动态策略灵活配置混合C'
$q1.Range('D29').Value = '0.13'
$q1.Range('E29').Value = '64.82'
$q1.Range('F29').Value = '2.78'
$q1.Range('G29').Value = '0.0036'
$q1.Cells.Item(29, 8).Value = 9
$q1.Cells.Item(30, 1).Value = 28
$q1.Range('B30').Value = '001608'
$q1.Range('C30').Value = '英大策略优选混合C'
$q1.Range('D30').Value = '0.03'
$q1.Range('E30').Value = '89.86'
$q1.Range('F30').Value = '7.21'
$q1.Range('G30').Value = '0.0022'
$q1.Cells.Item(30, 8).Value = 4
$q1.Cells.Item(31, 1).Value = 29
$q1.Range('B31').Value = '002035'
$q1.Range('C31').Value = '安信平稳增长混合C'
$q1.Range('D31').Value = '0.00'
$q1.Range('E31').Value = '65.16'
$q1.Range('F31').Value = '5.39'
$q1.Range('G31').NumberFormat = 'General'
$q1.Range('G31').Value = 0
$q1.Cells.Item(31, 8).Value = 4

# replicate the workbook's header / index-column style (bold, bordered,
# centered) instead of leaving these as freshly-synthesised styles
$styleSrcHeader.Copy()
$q1.Range('B1:H1').PasteSpecial(-4122)
$styleSrcIndex.Copy()
$q1.Range('A2:A31').PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- rebuilt '总计' summary sheet, inserted right after '2022-Q1' ----
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = '总计'

$total.Range('B1').Value = '日期'
$total.Range('C1').Value = '持有数量(只)'
$total.Range('D1').Value = '持有市值(亿元)'

# data rows
$total.Cells.Item(2, 1).Value = 0
$total.Range('B2').Value = '2022-Q1'
$total.Range('C2').Value = 30
$total.Range('D2').Value = 2.9
$total.Cells.Item(3, 1).Value = 1
$total.Range('B3').Value = '2021-Q4'
$total.Range('C3').Value = 50
$total.Range('D3').Value = 7.56
$total.Cells.Item(4, 1).Value = 2
$total.Range('B4').Value = '2021-Q3'
$total.Range('C4').Value = 62
$total.Range('D4').Value = 13.59
$total.Cells.Item(5, 1).Value = 3
$total.Range('B5').Value = '2021-Q2'
$total.Range('C5').Value = 41
$total.Range('D5').Value = 9.199999999999999
$total.Cells.Item(6, 1).Value = 4
$total.Range('B6').Value = '2021-Q1'
$total.Range('C6').Value = 64
$total.Range('D6').Value = 14.82
$total.Cells.Item(7, 1).Value = 5
$total.Range('B7').Value = '2020-Q4'
$total.Range('C7').Value = 34
$total.Range('D7').Value = 8.94

$styleSrcHeader.Copy()
$total.Range('B1:D1').PasteSpecial(-4122)
$styleSrcIndex.Copy()
$total.Range('A2:A7').PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q1.Range('A1').Select()
$total.Range('A1').Select()
